# Updated "Id" & changed the excel with display name
#
# Renames several header cells in row 1:
#   A1: "Product Id"  -> "Id"
#   C1: "Start Date"  -> "On-Site Date"
#   D1: "End Date"    -> "Off-Site Date"
#   H1: "Short Nm"    -> "Short Name"
# Then the columns are best-fit/auto-fit to the new content, and the
# active selection ends up on B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "On-Site Date"
$ws.Range("D1").Value = "Off-Site Date"
$ws.Range("H1").Value = "Short Name"

# Auto-fit (best-fit) the columns that changed width as a result of the
# new header/content ("Creation Date" column & "Display Name" column).
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("B1").Select() | Out-Null
